$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 4: add a new list paragraph after "...move up legs." paragraph
# ---------------------------------------------------------------------------
$p26 = $d.Paragraphs(26)
if ($p26.Range.Text -notmatch "move up legs") {
    throw "Paragraph 26 text mismatch: $($p26.Range.Text)"
}
$p26.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(27)
$newPara.Range.Text = "Verify all the previous steps for the three legs."

Write-Output "Change 4 done"

# ---------------------------------------------------------------------------
# Change 3: "The assets associated with the obstacles and scenery have
# already been designed." -- merge the split runs and underline "obstacles"
# ---------------------------------------------------------------------------
$p16 = $d.Paragraphs(16)
if ($p16.Range.Text -notmatch "The assets associated with the obstacles") {
    throw "Paragraph 16 text mismatch: $($p16.Range.Text)"
}
$p16start = $p16.Range.Start
$p16end = $p16.Range.End

# Locate the boundary right after "The asset" (but before the "s") so we can
# leave the leading tab+"The asset" run completely untouched (this preserves
# the <w:tab/> element, which gets flattened to a literal tab character
# whenever the run that owns it is edited).
$markerRng = $d.Range($p16start, $p16end)
if (-not $markerRng.Find.Execute("The asset", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    throw "Could not find 'The asset' in paragraph 16"
}
$afterAsset = $markerRng.End

# Merge everything after the leading tab+"The asset" run (keep the tab intact)
$rngBody = $d.Range($afterAsset, $p16end - 1)
if (-not $rngBody.Find.Execute("s associated with the obstacles and scenery have been already designed.", $true, $false, $false, $false, $false, $true, 1, $false, "s associated with the obstacles and scenery have been already designed.", 2)) {
    throw "Could not rebuild the body text of paragraph 16"
}

# underline the "obstacles" word that follows "associated with the "
$rngFind = $d.Range($afterAsset, $p16end - 1)
if (-not $rngFind.Find.Execute("associated with the obstacles", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    throw "Could not find 'obstacles' to underline in paragraph 16"
}
$obsStart = $rngFind.End - 9
$obsRng = $d.Range($obsStart, $rngFind.End)
if ($obsRng.Text -ne "obstacles") {
    throw "Underline target mismatch: $($obsRng.Text)"
}
$obsRng.Font.Underline = 1

Write-Output "Change 3 done"

# ---------------------------------------------------------------------------
# Change 2: "The user is in the level." -> "The user is in the level, and it
# has been already properly implemented."
# ---------------------------------------------------------------------------
$p14 = $d.Paragraphs(14)
if ($p14.Range.Text -notmatch "The user is in the level") {
    throw "Paragraph 14 text mismatch: $($p14.Range.Text)"
}
$p14start = $p14.Range.Start
$p14end = $p14.Range.End

# Only touch the trailing period run, leaving the tab+"The user is in the
# level" run untouched so the <w:tab/> element survives.
$markerRng2 = $d.Range($p14start, $p14end)
if (-not $markerRng2.Find.Execute("The user is in the level", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    throw "Could not find 'The user is in the level' in paragraph 14"
}
$periodRng = $d.Range($markerRng2.End, $p14end - 1)
if ($periodRng.Text -ne ".") {
    throw "Expected trailing period, found: $($periodRng.Text)"
}
$periodRng.Text = ", and it has been already properly implemented."

Write-Output "Change 2 done"

# ---------------------------------------------------------------------------
# Change 1: "...FR009 (Legs), FR004 (Obstacles)..." -> "...FR009 legs, FR004
# (Obstacles)..." with "Obstacles" now underlined
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(6)
if ($p6.Range.Text -notmatch "FR009") {
    throw "Paragraph 6 text mismatch: $($p6.Range.Text)"
}
$p6start = $p6.Range.Start
$p6end = $p6.Range.End

$rng1 = $d.Range($p6start, $p6end)
if (-not $rng1.Find.Execute("(Legs)", $true, $false, $false, $false, $false, $true, 1, $false, "legs", 2)) {
    throw "Could not find '(Legs)' in paragraph 6"
}

$rng2 = $d.Range($p6start, $p6end)
if (-not $rng2.Find.Execute("Obstacles", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    throw "Could not find 'Obstacles' to underline in paragraph 6"
}
if ($rng2.Text -ne "Obstacles") {
    throw "Underline target mismatch: $($rng2.Text)"
}
$rng2.Font.Underline = 1

Write-Output "Change 1 done"
